$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.802.63'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.066.40'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.19'
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.08'
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.063.51'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("E9").Value = '  -2.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.38'
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("E11").Value = '  -1.81%  '
$ws.Range("E12").Value = '  -3.10%  '
$ws.Range("E13").Value = '  -2.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.56'
$ws.Range("E14").Value = '  -4.30%  '
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.577.85'
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.716.21'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.80'
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.061.90'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '490.89'
$ws.Range("E21").Value = '  +3.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  -3.64%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.66'
$ws.Range("E24").Value = '  -1.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.67'
$ws.Range("E25").Value = '  -6.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.13'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.74'
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("E30").Value = '  -5.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.60'
$ws.Range("E31").Value = '  -2.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.45'
$ws.Range("E32").Value = '  -3.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.111'
$ws.Range("E33").Value = '  -3.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0908'
$ws.Range("E34").Value = '  -3.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.947'
$ws.Range("E36").Value = '  -2.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.58'
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '46.53'
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("E40").Value = '  -5.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.299'
$ws.Range("E41").Value = '  -3.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.31'
$ws.Range("E42").Value = '  -4.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.751.14'
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.91'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '365.94'
$ws.Range("E46").Value = '  -3.91%  '
$ws.Range("E47").Value = '  -4.69%  '
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.35'
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("E50").Value = '  -2.07%  '
$ws.Range("E51").Value = '  -1.96%  '
